# Applies the "po 04. 10. 2021" update to the Slovakia Covid daily-stats sheet:
#  - fills in previously-missing AgTests/AgPosit (F/G) values for rows 190 & 211
#  - corrects a batch of F/G figures for rows 268, 293, 518-574 (revised AG test counts)
#  - appends four new daily rows (575-578) for 2021-09-30 .. 2021-10-03
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections / fill-ins on existing rows ---
$ws.Cells.Item(190, 6).Value = 12
$ws.Cells.Item(190, 7).Value = 0
$ws.Cells.Item(211, 6).Value = 16
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(268, 6).Value = 17829
$ws.Cells.Item(293, 6).Value = 80183
$ws.Cells.Item(518, 6).Value = 7306
$ws.Cells.Item(519, 6).Value = 8100
$ws.Cells.Item(520, 6).Value = 10546
$ws.Cells.Item(521, 6).Value = 7013
$ws.Cells.Item(522, 6).Value = 5293
$ws.Cells.Item(523, 6).Value = 10369
$ws.Cells.Item(524, 6).Value = 7934
$ws.Cells.Item(525, 6).Value = 7741
$ws.Cells.Item(526, 6).Value = 8960
$ws.Cells.Item(527, 6).Value = 11718
$ws.Cells.Item(527, 7).Value = 34
$ws.Cells.Item(528, 6).Value = 8231
$ws.Cells.Item(529, 6).Value = 5879
$ws.Cells.Item(530, 6).Value = 12920
$ws.Cells.Item(531, 6).Value = 9374
$ws.Cells.Item(532, 6).Value = 10393
$ws.Cells.Item(533, 6).Value = 11956
$ws.Cells.Item(534, 6).Value = 16910
$ws.Cells.Item(535, 6).Value = 10246
$ws.Cells.Item(536, 6).Value = 8057
$ws.Cells.Item(537, 6).Value = 13833
$ws.Cells.Item(538, 6).Value = 11336
$ws.Cells.Item(539, 6).Value = 10676
$ws.Cells.Item(540, 6).Value = 12545
$ws.Cells.Item(541, 6).Value = 16707
$ws.Cells.Item(542, 6).Value = 10407
$ws.Cells.Item(542, 7).Value = 51
$ws.Cells.Item(543, 6).Value = 4763
$ws.Cells.Item(544, 6).Value = 14381
$ws.Cells.Item(545, 6).Value = 16697
$ws.Cells.Item(546, 6).Value = 3960
$ws.Cells.Item(547, 6).Value = 14018
$ws.Cells.Item(548, 6).Value = 17248
$ws.Cells.Item(549, 6).Value = 10788
$ws.Cells.Item(550, 6).Value = 8571
$ws.Cells.Item(550, 7).Value = 88
$ws.Cells.Item(551, 6).Value = 17639
$ws.Cells.Item(552, 6).Value = 15621
$ws.Cells.Item(553, 6).Value = 15339
$ws.Cells.Item(554, 6).Value = 17952
$ws.Cells.Item(554, 7).Value = 184
$ws.Cells.Item(555, 6).Value = 21490
$ws.Cells.Item(555, 7).Value = 184
$ws.Cells.Item(556, 6).Value = 12211
$ws.Cells.Item(557, 6).Value = 10962
$ws.Cells.Item(557, 7).Value = 149
$ws.Cells.Item(558, 6).Value = 24527
$ws.Cells.Item(559, 6).Value = 22481
$ws.Cells.Item(560, 6).Value = 6003
$ws.Cells.Item(561, 6).Value = 24049
$ws.Cells.Item(562, 6).Value = 27070
$ws.Cells.Item(562, 7).Value = 279
$ws.Cells.Item(563, 6).Value = 14092
$ws.Cells.Item(564, 6).Value = 14299
$ws.Cells.Item(565, 6).Value = 28601
$ws.Cells.Item(566, 6).Value = 25744
$ws.Cells.Item(567, 6).Value = 23385
$ws.Cells.Item(568, 6).Value = 23598
$ws.Cells.Item(568, 7).Value = 296
$ws.Cells.Item(569, 6).Value = 31991
$ws.Cells.Item(569, 7).Value = 360
$ws.Cells.Item(570, 6).Value = 15003
$ws.Cells.Item(570, 7).Value = 225
$ws.Cells.Item(571, 6).Value = 14984
$ws.Cells.Item(572, 6).Value = 32913
$ws.Cells.Item(572, 7).Value = 593
$ws.Cells.Item(573, 6).Value = 26514
$ws.Cells.Item(573, 7).Value = 399
$ws.Cells.Item(574, 6).Value = 23093
$ws.Cells.Item(574, 7).Value = 346

# --- New rows appended at the bottom ---
# Row 575
$ws.Cells.Item(575, 1).Value = 44469
$ws.Cells.Item(575, 2).Value = 413723
$ws.Cells.Item(575, 3).Value = 9308
$ws.Cells.Item(575, 4).Value = 1216
$ws.Cells.Item(575, 5).Value = 12649
$ws.Cells.Item(575, 6).Value = 24999
$ws.Cells.Item(575, 7).Value = 368
# Row 576
$ws.Cells.Item(576, 1).Value = 44470
$ws.Cells.Item(576, 2).Value = 415016
$ws.Cells.Item(576, 3).Value = 11203
$ws.Cells.Item(576, 4).Value = 1293
$ws.Cells.Item(576, 5).Value = 12660
$ws.Cells.Item(576, 6).Value = 25209
$ws.Cells.Item(576, 7).Value = 383
# Row 577
$ws.Cells.Item(577, 1).Value = 44471
$ws.Cells.Item(577, 2).Value = 415993
$ws.Cells.Item(577, 3).Value = 7308
$ws.Cells.Item(577, 4).Value = 977
$ws.Cells.Item(577, 5).Value = 12668
$ws.Cells.Item(577, 6).Value = 11672
$ws.Cells.Item(577, 7).Value = 221
# Row 578
$ws.Cells.Item(578, 1).Value = 44472
$ws.Cells.Item(578, 2).Value = 416260
$ws.Cells.Item(578, 3).Value = 2223
$ws.Cells.Item(578, 4).Value = 267
$ws.Cells.Item(578, 5).Value = 12676
$ws.Cells.Item(578, 6).Value = 8708
$ws.Cells.Item(578, 7).Value = 206
